# Apply updated cryptocurrency price/volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) are stored as text in the sheet (not numbers).
# Temporarily force column D to Text format so numeric-looking values (e.g. "255.40")
# are written as strings instead of being auto-converted to floating point numbers,
# then restore the original (default) style so the cell formatting matches the source.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '35.300.22'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '1.911.76'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('E5').Value = '  +9.16%  '
$ws.Range('D6').Value = '255.40'
$ws.Range('E6').Value = '  +3.63%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +1.15%  '
$ws.Range('D9').Value = '0.368'
$ws.Range('E9').Value = '  +5.93%  '
$ws.Range('D10').Value = '53.26'
$ws.Range('E10').Value = '  -0.30%  '
$ws.Range('D11').Value = '0.0765'
$ws.Range('E11').Value = '  +6.41%  '
$ws.Range('D12').Value = '0.0987'
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('D13').Value = '13.16'
$ws.Range('E13').Value = '  +6.63%  '
$ws.Range('D14').Value = '2.191.29'
$ws.Range('E14').Value = '  +0.02%  '
$ws.Range('D15').Value = '0.739'
$ws.Range('E15').Value = '  +5.24%  '
$ws.Range('E16').Value = '  +4.07%  '
$ws.Range('D17').Value = '1.918.05'
$ws.Range('E17').Value = '  +0.27%  '
$ws.Range('D18').Value = '35.297.98'
$ws.Range('E18').Value = '  -0.23%  '
$ws.Range('D19').Value = '75.17'
$ws.Range('E19').Value = '  +4.27%  '
$ws.Range('D20').Value = '0.0₃0849'
$ws.Range('E20').Value = '  +3.52%  '
$ws.Range('D21').Value = '245.90'
$ws.Range('E21').Value = '  +1.83%  '
$ws.Range('E22').Value = '  +4.96%  '
$ws.Range('E23').Value = '  +6.92%  '
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('E25').Value = '  +7.40%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').Value = '166.77'
$ws.Range('E27').Value = '  -2.41%  '
$ws.Range('E28').Value = '  +4.08%  '
$ws.Range('E29').Value = '  +2.23%  '
$ws.Range('D30').Value = '0.133'
$ws.Range('E30').Value = '  +4.30%  '
$ws.Range('D31').Value = '4.128.98'
$ws.Range('E31').Value = '  -0.67%  '
$ws.Range('E32').Value = '  +27.53%  '
$ws.Range('D33').Value = '4.37'
$ws.Range('E33').Value = '  +5.20%  '
$ws.Range('E34').Value = '  +14.87%  '
$ws.Range('D35').Value = '0.0593'
$ws.Range('E35').Value = '  +4.62%  '
$ws.Range('D36').Value = '4.28'
$ws.Range('E36').Value = '  +4.27%  '
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('E38').Value = '  -3.54%  '
$ws.Range('E39').Value = '  +0.10%  '
$ws.Range('D40').Value = '100.35'
$ws.Range('E40').Value = '  +11.56%  '
$ws.Range('E41').Value = '  +6.38%  '
$ws.Range('D42').Value = '17.09'
$ws.Range('E42').Value = '  +5.35%  '
$ws.Range('E43').Value = '  +1.98%  '
$ws.Range('E44').Value = '  +0.25%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = '2.47'
$ws.Range('E45').Value = '  +3.16%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '1.345.20'
$ws.Range('E46').Value = '  +0.41%  '
$ws.Range('E47').Value = '  +1.32%  '
$ws.Range('E48').Value = '  +3.55%  '
$ws.Range('E49').Value = '  -0.78%  '
$ws.Range('D50').Value = '45.02'
$ws.Range('E50').Value = '  -8.43%  '
$ws.Range('D51').Value = '0.0759'
$ws.Range('E51').Value = '  +7.21%  '

$priceRange.Style = "Normal"

